$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 412.66666
$ws.Range("I6").Value = 401
$ws.Range("J6").Value = 459.33334
$ws.Range("K6").Value = 1203
$ws.Range("L6").Value = 1378.00002
$ws.Range("M6").Value = -1091
$ws.Range("N6").Value = -1602.00002

# Sheet ALC, Row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 985.9
$ws.Range("I28").Value = 985.9
$ws.Range("K28").Value = 985.9
$ws.Range("M28").Value = -500.9

# Sheet ALC, Row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 103.916664
$ws.Range("I33").Value = 99.72727
$ws.Range("K33").Value = 99.72727
$ws.Range("M33").Value = 129.27273

# Sheet ALC, Row 34
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2579.4
$ws.Range("I34").Value = 1974.25
$ws.Range("K34").Value = 1974.25
$ws.Range("M34").Value = -1771.25

# Sheet ALC, Row 36
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 2579.4
$ws.Range("I36").Value = 1974.25
$ws.Range("K36").Value = 1974.25
$ws.Range("M36").Value = -1259.25

# Sheet ALC, Row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 334.22223
$ws.Range("I38").Value = 126
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 378
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -6
$ws.Range("N38").Value = -6744

# Sheet ALC, Row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2113.9167
$ws.Range("I40").Value = 1991.4445
$ws.Range("J40").Value = 2481.3333
$ws.Range("K40").Value = 1991.4445
$ws.Range("L40").Value = 2481.3333
$ws.Range("M40").Value = -1816.4445
$ws.Range("N40").Value = -2831.3333

# Sheet ALC, Row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 683
$ws.Range("I92").Value = 683
$ws.Range("K92").Value = 683
$ws.Range("M92").Value = 565

# Sheet ALC, Row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1092
$ws.Range("I32").Value = 1092
$ws.Range("K32").Value = 1092
$ws.Range("M32").Value = -805

# Sheet ARM, Row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3785.8
$ws.Range("I35").Value = 3732.5
$ws.Range("K35").Value = 3732.5
$ws.Range("M35").Value = -3326.5

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4678
$ws.Range("I61").Value = 4678
$ws.Range("K61").Value = 4678
$ws.Range("M61").Value = -4466

# Sheet ARM, Row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 205767.67
$ws.Range("J101").Value = 205767.67
$ws.Range("L101").Value = 205767.67
$ws.Range("N101").Value = -212257.67

# Sheet ARM, Row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4236.1665
$ws.Range("I102").Value = 4236.1665
$ws.Range("K102").Value = 4236.1665
$ws.Range("M102").Value = -2614.1665

# Sheet ARM, Row 127
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H127").Value = 36000
$ws.Range("I127").Value = 36000
$ws.Range("K127").Value = 36000
$ws.Range("M127").Value = -31040

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3119.6
$ws.Range("I132").Value = 3119.6
$ws.Range("K132").Value = 9358.799999999999
$ws.Range("M132").Value = -6828.799999999999

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4678
$ws.Range("I136").Value = 4678
$ws.Range("K136").Value = 14034
$ws.Range("M136").Value = -11484

# Sheet BSM, Row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1649.4445
$ws.Range("I105").Value = 1705.625
$ws.Range("K105").Value = 1705.625
$ws.Range("M105").Value = 41.375

# Sheet BSM, Row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9340.299999999999
$ws.Range("I134").Value = 7300.375
$ws.Range("K134").Value = 21901.125
$ws.Range("M134").Value = -19366.125

# Sheet CRP, Row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 835.6667
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# Sheet CRP, Row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# Sheet CRP, Row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 277.66666
$ws.Range("I25").Value = 200
$ws.Range("K25").Value = 200
$ws.Range("M25").Value = -26

# Sheet CRP, Row 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 14685.143
$ws.Range("I39").Value = 4559.2
$ws.Range("K39").Value = 4559.2
$ws.Range("M39").Value = -4168.2

# Sheet CRP, Row 49
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 14685.143
$ws.Range("I49").Value = 4559.2
$ws.Range("K49").Value = 4559.2
$ws.Range("M49").Value = -4377.2

# Sheet CRP, Row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7333
$ws.Range("J86").Value = 7500
$ws.Range("L86").Value = 7500
$ws.Range("N86").Value = -9746

# Sheet CRP, Row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 22190.334
$ws.Range("J88").Value = 22190.334
$ws.Range("L88").Value = 22190.334
$ws.Range("N88").Value = -23002.334

# Sheet CRP, Row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7333
$ws.Range("J89").Value = 7500
$ws.Range("L89").Value = 37500
$ws.Range("N89").Value = -48732

# Sheet CRP, Row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 22190.334
$ws.Range("J91").Value = 22190.334
$ws.Range("L91").Value = 22190.334
$ws.Range("N91").Value = -24998.334

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2298.875
$ws.Range("I132").Value = 2046.1666
$ws.Range("J132").Value = 3057
$ws.Range("K132").Value = 6138.4998
$ws.Range("L132").Value = 9171
$ws.Range("M132").Value = -3608.4998
$ws.Range("N132").Value = -14231

# Sheet CUL, Row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 327.2
$ws.Range("I5").Value = 322.75
$ws.Range("J5").Value = 345
$ws.Range("K5").Value = 968.25
$ws.Range("L5").Value = 1035
$ws.Range("M5").Value = -856.25
$ws.Range("N5").Value = -1259

# Sheet CUL, Row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 183.2
$ws.Range("J12").Value = 228.5
$ws.Range("L12").Value = 685.5
$ws.Range("N12").Value = -1031.5

# Sheet CUL, Row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2843.375
$ws.Range("J34").Value = 4749.75
$ws.Range("L34").Value = 14249.25
$ws.Range("N34").Value = -14417.25

# Sheet CUL, Row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2716.25
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 3338.3333
$ws.Range("K46").Value = 2550
$ws.Range("L46").Value = 10014.9999
$ws.Range("M46").Value = -2459
$ws.Range("N46").Value = -10196.9999

# Sheet CUL, Row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 259.66666
$ws.Range("J122").Value = 362.5
$ws.Range("L122").Value = 3262.5
$ws.Range("N122").Value = -8162.5

# Sheet CUL, Row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 327.2
$ws.Range("I135").Value = 322.75
$ws.Range("J135").Value = 345
$ws.Range("K135").Value = 2904.75
$ws.Range("L135").Value = 3105
$ws.Range("M135").Value = -369.75
$ws.Range("N135").Value = -8175

# Sheet GSM, Row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 160.75
$ws.Range("I2").Value = 176.57143
$ws.Range("K2").Value = 176.57143
$ws.Range("M2").Value = -63.57142999999999

# Sheet GSM, Row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4816924
$ws.Range("J11").Value = 410000
$ws.Range("L11").Value = 410000
$ws.Range("N11").Value = -410278

# Sheet GSM, Row 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 409
$ws.Range("I31").Value = 409
$ws.Range("K31").Value = 409
$ws.Range("M31").Value = -117

# Sheet GSM, Row 37
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 409
$ws.Range("I37").Value = 409
$ws.Range("K37").Value = 409
$ws.Range("M37").Value = -132

# Sheet GSM, Row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4714.7
$ws.Range("I80").Value = 4057.8333
$ws.Range("J80").Value = 5700
$ws.Range("K80").Value = 4057.8333
$ws.Range("L80").Value = 5700
$ws.Range("M80").Value = -3059.8333
$ws.Range("N80").Value = -7696

# Sheet GSM, Row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4714.7
$ws.Range("I83").Value = 4057.8333
$ws.Range("J83").Value = 5700
$ws.Range("K83").Value = 20289.1665
$ws.Range("L83").Value = 28500
$ws.Range("M83").Value = -15297.1665
$ws.Range("N83").Value = -38484

# Sheet GSM, Row 101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 24500
$ws.Range("J101").Value = 24500
$ws.Range("L101").Value = 24500
$ws.Range("N101").Value = -30990

# Sheet GSM, Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5002.3335
$ws.Range("I122").Value = 5002.3335
$ws.Range("K122").Value = 15007.0005
$ws.Range("M122").Value = -12557.0005

# Sheet GSM, Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1046.5
$ws.Range("I132").Value = 1046.5
$ws.Range("K132").Value = 3139.5
$ws.Range("M132").Value = -609.5

# Sheet LTW, Row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 2015000
$ws.Range("J24").Value = 2015000
$ws.Range("L24").Value = 2015000
$ws.Range("N24").Value = -2015686

# Sheet LTW, Row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 558.1667
$ws.Range("I55").Value = 656.6667
$ws.Range("J55").Value = 459.66666
$ws.Range("K55").Value = 656.6667
$ws.Range("L55").Value = 459.66666
$ws.Range("M55").Value = -483.6667
$ws.Range("N55").Value = -805.66666

# Sheet LTW, Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3332.7778
$ws.Range("I122").Value = 3249.375
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 9748.125
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7298.125
$ws.Range("N122").Value = -16900

# Sheet LTW, Row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 999
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 2997
$ws.Range("M132").Value = -467

# Sheet WVR, Row 97
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 60000
$ws.Range("J97").Value = 60000
$ws.Range("L97").Value = 60000
$ws.Range("N97").Value = -61982

# Sheet WVR, Row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 350
$ws.Range("I122").Value = 350
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1050
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1400
$ws.Range("N122").ClearContents()

# Sheet WVR, Row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2179.1428
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1939.8
$ws.Range("I132").Value = 1424.75
$ws.Range("K132").Value = 4274.25
$ws.Range("M132").Value = -1744.25

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13888.111
$ws.Range("J136").Value = 14999.667
$ws.Range("L136").Value = 44999.001
$ws.Range("N136").Value = -50099.001
